# Updated cases for new parameters
#
# andes/cases/ieee14/ieee14_solar.xlsx
#
# 1) REGCA1 sheet: the "Iqmax" / "Iqmin" parameter columns are removed
#    (they were never used by the REGCA1 model), shifting "ra"/"xs" left.
# 2) REPCA1 (the sheet that holds the REECA1-style renewable-electrical-
#    control parameters) sheet: a new "PLflag" parameter column is
#    inserted right after "Fflag" (and before "Tfltr"), with a value of 0
#    for the existing device row.

$wb = $excel.ActiveWorkbook

# --- 1) REGCA1: drop the Iqmax / Iqmin columns (U:V) ---
$wsRegca1 = $wb.Worksheets.Item("REGCA1")
$wsRegca1.Range("U1:V1").EntireColumn.Delete()

# --- 2) REPCA1: insert the new "PLflag" column before the old column L (Tfltr) ---
$wsRepca1 = $wb.Worksheets.Item("REPCA1")
$wsRepca1.Range("L1").EntireColumn.Insert()
$wsRepca1.Range("L1").Value = "PLflag"
$wsRepca1.Range("L2").Value = 0

# --- refresh the view/selection state to match the saved workbook ---
$wsRegca1.Range("R7").Select()

$wsRepca1.Activate()
$wsRepca1.Range("L2").Select()
